$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay stored as text (matching the
# original inlineStr cells). Force text format before writing, then restore the
# default ("Normal") cell style so no stray style index is left behind.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D13", "D14", "D19", "D20", "D21", "D22", "D24", "D28", "D29", "D30", "D31", "D32", "D33", "D40", "D41", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.642.80"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.298.93"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "312.38"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "104.46"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  -2.50%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "39.79"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "0.985"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "2.648.01"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "2.300.10"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "42.626.19"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "7.31"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "13.60"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0000104"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("E23").Value = "  -5.29%  "
$ws.Range("D24").Value = "267.27"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "7.14"
$ws.Range("E28").Value = "  +15.36%  "
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "22.34"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").Value = "36.12"
$ws.Range("E31").Value = "  -5.17%  "
$ws.Range("D32").Value = "164.79"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "0.0853"
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").Value = "3.65"
$ws.Range("E40").Value = "  -2.75%  "
$ws.Range("D41").Value = "107.61"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "71.14"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "12.16"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("D47").Value = "1.740.83"
$ws.Range("E47").Value = "  +9.30%  "
$ws.Range("D48").Value = "110.59"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("D49").Value = "77.69"
$ws.Range("E49").Value = "  -6.39%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.64"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "5.15"
$ws.Range("E51").Value = "  -3.24%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
